$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert a new row 3 on sheet 1 (performance comparison sheet) for the
# "ID3" algorithm, pushing the existing C4.5 / C5.0 rows down to rows 4/5.
# Inserting via Rows.Insert() also copies the row-2 formatting (style "2")
# down onto the new row automatically, matching the target cell styles.
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "ID3"
$ws1.Range("B3").Value = 0.96640000000000004
$ws1.Range("C3").Value = 0.77559944000000003
$ws1.Range("D3").Value = 0.78080000000000005
$ws1.Range("E3").Value = 0.77736518589541104
$ws1.Range("F3").Value = 0.55412151216565597
$ws1.Range("G3").Value = 0.533917894411294

# Make the first worksheet the active one (it was previously the second
# sheet that was active/selected) and restore the saved selection/cursor
# position on it. Sheet 2 keeps its original selection, it just stops
# being the active tab as a side effect of activating sheet 1.
$ws1.Activate()
$ws1.Range("D18").Select()
